# AddEventSD.pptx - "Implementation of Calendar and event management" diagram.
#
# Two shapes on slide 1 are adjusted:
#  - "Straight Arrow Connector 49" (shape id 50): the connector is
#    repositioned/resized and its vertical flip is cleared so it is
#    rendered as a (near) horizontal arrow instead of the old
#    diagonal/flipped one.
#  - "TextBox 72" (shape id 73, the "addEvent()" label): moved down so it
#    sits correctly against the updated connector.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# "Straight Arrow Connector 49"
$connector = Get-ShapeById $s 50
$connector.Left   = 141.8589713780
$connector.Top    = 232.8859025118
$connector.Width  = 148.7726831654
$connector.Height = 0.0000787402
$connector.VerticalFlip = 0

# "TextBox 72" ("addEvent()" label)
$label = Get-ShapeById $s 73
$label.Left = 312.1407016614
$label.Top  = 222.3104784409
